$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2
$ws.Range("F2").Value = 29
$ws.Range("H2").Value = 41

# Row 18
$ws.Range("F18").Value = 46
$ws.Range("H18").Value = 82

# Row 26
$ws.Range("E26").Value = 26
$ws.Range("F26").Value = 12
$ws.Range("H26").Value = 22

# Row 28
$ws.Range("F28").Value = 13
$ws.Range("H28").Value = 15

# Row 38
$ws.Range("E38").Value = 73

# Row 48
$ws.Range("E48").Value = 29
$ws.Range("F48").Value = 20
$ws.Range("H48").Value = 25

# Row 51
$ws.Range("F51").Value = 9
$ws.Range("H51").Value = 9

# Row 71
$ws.Range("F71").Value = 17
$ws.Range("H71").Value = 27

# Row 78
$ws.Range("F78").Value = 20
$ws.Range("H78").Value = 41

# Row 81
$ws.Range("F81").Value = 8
$ws.Range("H81").Value = 13

# Row 89
$ws.Range("F89").Value = 15
$ws.Range("H89").Value = 22
